$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 22623.215
$ws.Range("J87").Value = 22623.215
$ws.Range("L87").Value = 22623.215
$ws.Range("N87").Value = -25119.215
$ws.Range("H90").Value = 22623.215
$ws.Range("J90").Value = 22623.215
$ws.Range("L90").Value = 67869.645
$ws.Range("N90").Value = -80349.645
$ws.Range("H113").Value = 7993.3335
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 8792
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 8792
$ws.Range("M113").Value = -746
$ws.Range("N113").Value = -15300
$ws.Range("H116").Value = 1011600.5
$ws.Range("I116").Value = 5001102.5
$ws.Range("J116").Value = 14225
$ws.Range("K116").Value = 5001102.5
$ws.Range("L116").Value = 14225
$ws.Range("M116").Value = -4997660.5
$ws.Range("N116").Value = -21109
$ws.Range("H132").Value = 35720092
$ws.Range("I132").Value = 40005504
$ws.Range("J132").Value = 8333.333000000001
$ws.Range("K132").Value = 120016512
$ws.Range("L132").Value = 24999.999
$ws.Range("M132").Value = -120013982
$ws.Range("N132").Value = -30059.999
$ws.Range("H137").Value = 1702938.8
$ws.Range("J137").Value = 4443.6665
$ws.Range("L137").Value = 13330.9995
$ws.Range("N137").Value = -18430.9995

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 625
$ws.Range("I2").Value = 625
$ws.Range("K2").Value = 625
$ws.Range("M2").Value = -512
$ws.Range("H61").Value = 2196.9167
$ws.Range("I61").Value = 1023.75
$ws.Range("J61").Value = 2783.5
$ws.Range("K61").Value = 1023.75
$ws.Range("L61").Value = 2783.5
$ws.Range("M61").Value = -811.75
$ws.Range("N61").Value = -3207.5
$ws.Range("H110").Value = 1703.6666
$ws.Range("I110").Value = 1703.6666
$ws.Range("K110").Value = 1703.6666
$ws.Range("M110").Value = 341.3334
$ws.Range("H116").Value = 625
$ws.Range("I116").Value = 625
$ws.Range("K116").Value = 625
$ws.Range("M116").Value = 1669
$ws.Range("H132").Value = 2342.8262
$ws.Range("I132").Value = 1067.5
$ws.Range("J132").Value = 3734.0908
$ws.Range("K132").Value = 3202.5
$ws.Range("L132").Value = 11202.2724
$ws.Range("M132").Value = -672.5
$ws.Range("N132").Value = -16262.2724
$ws.Range("H133").Value = 28926.666
$ws.Range("J133").Value = 28926.666
$ws.Range("L133").Value = 28926.666
$ws.Range("N133").Value = -33986.666
$ws.Range("H136").Value = 2196.9167
$ws.Range("I136").Value = 1023.75
$ws.Range("J136").Value = 2783.5
$ws.Range("K136").Value = 3071.25
$ws.Range("L136").Value = 8350.5
$ws.Range("M136").Value = -521.25
$ws.Range("N136").Value = -13450.5
$ws.Range("H137").Value = 39890
$ws.Range("J137").Value = 39890
$ws.Range("L137").Value = 39890
$ws.Range("N137").Value = -50090
$ws.Range("H139").Value = 43920.6
$ws.Range("J139").Value = 43920.6
$ws.Range("L139").Value = 43920.6
$ws.Range("N139").Value = -54200.6

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 625
$ws.Range("I3").Value = 625
$ws.Range("K3").Value = 625
$ws.Range("M3").Value = -511
$ws.Range("H107").Value = 789.28
$ws.Range("I107").Value = 810.087
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 810.087
$ws.Range("L107").Value = 550
$ws.Range("M107").Value = 1109.913
$ws.Range("N107").Value = -4390
$ws.Range("H137").Value = 45780
$ws.Range("J137").Value = 45780
$ws.Range("L137").Value = 45780
$ws.Range("N137").Value = -55980

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3181.5217
$ws.Range("I31").Value = 1072.5625
$ws.Range("K31").Value = 1072.5625
$ws.Range("M31").Value = -777.5625
$ws.Range("H34").Value = 3181.5217
$ws.Range("I34").Value = 1072.5625
$ws.Range("K34").Value = 1072.5625
$ws.Range("M34").Value = -870.5625
$ws.Range("H107").Value = 791.9167
$ws.Range("I107").Value = 687.7778
$ws.Range("J107").Value = 1104.3334
$ws.Range("K107").Value = 687.7778
$ws.Range("L107").Value = 1104.3334
$ws.Range("M107").Value = 1232.2222
$ws.Range("N107").Value = -4944.3334
$ws.Range("H132").Value = 2531.8635
$ws.Range("I132").Value = 1575.1875
$ws.Range("J132").Value = 5083
$ws.Range("K132").Value = 4725.5625
$ws.Range("L132").Value = 15249
$ws.Range("M132").Value = -2195.5625
$ws.Range("N132").Value = -20309
$ws.Range("H137").Value = 46728.332
$ws.Range("J137").Value = 46728.332
$ws.Range("L137").Value = 46728.332
$ws.Range("N137").Value = -56928.332

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 581380.75
$ws.Range("J5").Value = 835449.9
$ws.Range("L5").Value = 2506349.7
$ws.Range("N5").Value = -2506573.7
$ws.Range("H131").Value = 746.1
$ws.Range("J131").Value = 797.8876299999999
$ws.Range("L131").Value = 2393.66289
$ws.Range("N131").Value = -12473.66289
$ws.Range("H135").Value = 581380.75
$ws.Range("J135").Value = 835449.9
$ws.Range("L135").Value = 7519049.100000001
$ws.Range("N135").Value = -7524119.100000001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 22730064
$ws.Range("I80").Value = 31252714
$ws.Range("J80").Value = 2999.6667
$ws.Range("K80").Value = 31252714
$ws.Range("L80").Value = 2999.6667
$ws.Range("M80").Value = -31251716
$ws.Range("N80").Value = -4995.6667
$ws.Range("H83").Value = 22730064
$ws.Range("I83").Value = 31252714
$ws.Range("J83").Value = 2999.6667
$ws.Range("K83").Value = 156263570
$ws.Range("L83").Value = 14998.3335
$ws.Range("M83").Value = -156258578
$ws.Range("N83").Value = -24982.3335
$ws.Range("H102").Value = 3670.3
$ws.Range("I102").Value = 2100.4285
$ws.Range("J102").Value = 7333.3335
$ws.Range("K102").Value = 2100.4285
$ws.Range("L102").Value = 7333.3335
$ws.Range("M102").Value = -478.4285
$ws.Range("N102").Value = -10577.3335
$ws.Range("H137").Value = 37155
$ws.Range("J137").Value = 37155
$ws.Range("L137").Value = 37155
$ws.Range("N137").Value = -47355

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4059.25
$ws.Range("I7").Value = 3352.7273
$ws.Range("J7").Value = 4922.778
$ws.Range("K7").Value = 3352.7273
$ws.Range("L7").Value = 4922.778
$ws.Range("M7").Value = -3240.7273
$ws.Range("N7").Value = -5146.778
$ws.Range("H40").Value = 5499.5
$ws.Range("I40").Value = 4999.375
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 4999.375
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -4863.375
$ws.Range("N40").Value = -7772
$ws.Range("H126").Value = 4059.25
$ws.Range("I126").Value = 3352.7273
$ws.Range("J126").Value = 4922.778
$ws.Range("K126").Value = 10058.1819
$ws.Range("L126").Value = 14768.334
$ws.Range("M126").Value = -7588.1819
$ws.Range("N126").Value = -19708.334
$ws.Range("H133").Value = 29022.75
$ws.Range("J133").Value = 29022.75
$ws.Range("L133").Value = 29022.75
$ws.Range("N133").Value = -34082.75
